$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Range("B2").Font.Size = 22
$ws2.Range("B4:H4").HorizontalAlignment = -4108
$ws2.Range("B4:H4").HorizontalAlignment = -4108
